# Autogenerated data refresh for "MSME Country Indicators - Portugal Summary"
# Updates a handful of percentage figures on the Summary sheet with more
# precise (2-decimal) values from the upstream source.
#
# These cells hold numeric-looking figures that are stored as TEXT in the
# workbook (e.g. "107.6"), so we use the classic Excel leading-apostrophe
# trick to keep them as text instead of letting them be auto-converted to
# numbers when we write the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people) — Source Type: Statistical Institution
$ws.Range("D11").Value = "'107.56"

# Employment (% of total) — Source Type: Statistical Institution
$ws.Range("D12").Value = "'77.35"

# Enterprises density (per 1000 people) — Source Type: SME Associations (Micro / SMEs / MSMEs)
$ws.Range("B33").Value = "'74.83"
$ws.Range("C33").Value = "'3.87"

# Employment (% of total) — Source Type: SME Associations (Micro / SMEs / MSMEs)
$ws.Range("B34").Value = "'41.39"
$ws.Range("C34").Value = "'37.64"
$ws.Range("D34").Value = "'79.03"

# Enterprises (% of total) — Source Type: SME Associations (Micro / SMEs / MSMEs)
$ws.Range("B36").Value = "'94.98"
$ws.Range("C36").Value = "'4.92"

# Value added to the economy (% of total) — Source Type: SME Associations
$ws.Range("B40").Value = "'23.08"
$ws.Range("C40").Value = "'44.65"
$ws.Range("D40").Value = "'67.73"
